# Select Random Player and fix budget calculation after each bid
$wb = $excel.ActiveWorkbook

# --- Bids sheet: the "Jasprit Bumrah" bid rows were replaced by a newly
#     (randomly) selected player, "Rohit Sharma". The old Virat Kohli
#     bid rows shift up to rows 2-3, and the new Rohit Sharma bids land
#     in rows 4-5 with their own amounts.
$bids = $wb.Worksheets.Item("Bids")
$bids.Range("A2").Value = "Virat Kohli"
$bids.Range("B2").Value = "Aniket"
$bids.Range("C2").Value = 10.0

$bids.Range("A3").Value = "Virat Kohli"
$bids.Range("B3").Value = "Hiren"
$bids.Range("C3").Value = 20.0

$bids.Range("A4").Value = "Rohit Sharma"
$bids.Range("B4").Value = "Anurag"
$bids.Range("C4").Value = 15.0

$bids.Range("A5").Value = "Rohit Sharma"
$bids.Range("B5").Value = "Vipul"
$bids.Range("C5").Value = 10.0

# --- Participants sheet: fix Anurag's remaining budget after his bid
#     on Rohit Sharma (100 base budget - 15 bid = 85).
$participants = $wb.Worksheets.Item("Participants")
$participants.Range("B4").Value = 85.0
